$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.422.64'
$ws.Range('E2').Value = '  +1.41%  '

$ws.Range('D3').Value = '3.916.39'
$ws.Range('E3').Value = '  -0.91%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '486.52'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.95%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.82'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.30%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.623'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.08%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.733'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.22%  '

$ws.Range('E10').Value = '  -0.58%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000346'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.16%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '43.23'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.54%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.77'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.85%  '

$ws.Range('D14').Value = '4.539.04'
$ws.Range('E14').Value = '  -0.82%  '

$ws.Range('D15').Value = '3.912.49'
$ws.Range('E15').Value = '  -1.22%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.36'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.42%  '

$ws.Range('E17').Value = '  -1.25%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '20.05'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.24%  '

$ws.Range('E19').Value = '  -1.33%  '

$ws.Range('D20').Value = '68.468.39'
$ws.Range('E20').Value = '  +1.25%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '434.45'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.00%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.49'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.20%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '15.15'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +4.03%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '88.28'
$ws.Range('D24').ClearFormats()

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.48'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +16.49%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.19'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +10.26%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.65'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.20%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '37.97'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.86%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.75'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.22%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '721.91'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.42%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '13.81'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.34%  '

$ws.Range('E32').Value = '  -2.58%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.93'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.83%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.24'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +15.76%  '

$ws.Range('B35').Value = 'PEPE'
$ws.Range('C35').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D35').Value = '0.0₃0882'
$ws.Range('E35').Value = '  +3.30%  '

$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '41.42'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.22%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '60.89'
$ws.Range('D37').ClearFormats()

$ws.Range('E38').Value = '  -3.54%  '

$ws.Range('E39').Value = '  +0.14%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.394'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +16.67%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0493'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.39%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.96'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +16.73%  '

$ws.Range('E43').Value = '  +1.24%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.99'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +5.30%  '

$ws.Range('E45').Value = '  +5.20%  '

$ws.Range('E46').Value = '  -2.11%  '

$ws.Range('E47').Value = '  +0.08%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.42'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.40%  '

$ws.Range('E49').Value = '  -4.45%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '145.31'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.16%  '

$ws.Range('E51').Value = '  +28.64%  '
